$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Z2:Z7").Value = "2025-11-02T02:03:48.251009"
$ws.Range("Z8:Z17").Value = "2025-11-02T02:03:48.252010"
$ws.Range("Z18:Z27").Value = "2025-11-02T02:03:48.253009"
$ws.Range("Z28:Z29").Value = "2025-11-02T02:03:48.254009"
